$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 37957
$ws1.Range("F5").Value = 795
$ws1.Range("F6").Value = 491
$ws1.Range("F9").Value = 866
$ws1.Range("F10").Value = 108
$ws1.Range("F11").Value = 750
$ws1.Range("F12").Value = 587
$ws1.Range("F13").Value = 82
$ws1.Range("F14").Value = 39
$ws1.Range("F15").Value = 37
$ws1.Range("F16").Value = 687
$ws1.Range("F18").Value = 489
$ws1.Range("F19").Value = 450
$ws1.Range("F22").Value = 879
$ws1.Range("F23").Value = 2594
$ws1.Range("F24").Value = 1076
$ws1.Range("F25").Value = 584
$ws1.Range("F26").Value = 121
$ws1.Range("F27").Value = 1180
$ws1.Range("F29").Value = 837
$ws1.Range("F30").Value = 78
$ws1.Range("F31").Value = 1183
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F3").Value = 444
$ws2.Range("F4").Value = 338
$ws2.Range("F10").Value = 17
$ws2.Range("F12").Value = 14
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 668
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 668
$ws4.Range("F3").Value = 37957
$ws4.Range("F6").Value = 795
$ws4.Range("F7").Value = 491
$ws4.Range("F11").Value = 444
$ws4.Range("F12").Value = 444
$ws4.Range("F16").Value = 866
$ws4.Range("F17").Value = 108
$ws4.Range("F18").Value = 750
$ws4.Range("F19").Value = 587
$ws4.Range("F20").Value = 82
$ws4.Range("F22").Value = 39
$ws4.Range("F25").Value = 17
$ws4.Range("F26").Value = 37
$ws4.Range("F28").Value = 687
$ws4.Range("F30").Value = 489
$ws4.Range("F31").Value = 450
$ws4.Range("F34").Value = 879
$ws4.Range("F35").Value = 2594
$ws4.Range("F36").Value = 1077
$ws4.Range("F37").Value = 584
$ws4.Range("F38").Value = 121
$ws4.Range("F39").Value = 1180
$ws4.Range("F41").Value = 14
$ws4.Range("F42").Value = 837
$ws4.Range("F43").Value = 78
$ws4.Range("F44").Value = 1183
